$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "56.338.94"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "3.008.15"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.13"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.06"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.13"
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.369"
$ws.Range("E11").Value = "  +3.53%  "
$ws.Range("D12").Value = "3.519.97"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.50"
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "56.318.59"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "3.005.54"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.97"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.28"
$ws.Range("E21").Value = "  +4.22%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.86"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("D25").Value = "3.131.67"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  +1.60%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.0₃0938"
$ws.Range("E28").Value = "  +5.27%  "
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.88"
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.41"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.86"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.58"
$ws.Range("E36").Value = "  +9.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.84"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.23"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "3.047.13"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.40"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").Value = "2.203.27"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0240"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.926"
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.85"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.49"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0852"
$ws.Range("E51").Value = "  -1.76%  "
